# Apply the authored change:
#  - Update the last data row (row 8) values for nufft512 (B8:E8)
#  - Move the sheet selection to reflect the cells the author last touched
#    (Excel records a multi-area selection here: E8 and B10 with B10 as the
#    active cell; this runtime's selection model only tracks a single
#    area/active-cell, so we land the selection on B10 - the final active
#    cell - as the closest reachable equivalent).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 6
$ws.Range("C8").Value = 9
$ws.Range("D8").Value = 9
$ws.Range("E8").Value = 16

$ws.Range("B10").Select()
